$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 313, pushing existing rows 313:341 down to 314:342
$ws.Rows("313:313").Insert()

# Populate the newly inserted row 313 with the new data record
$ws.Cells.Item(313, 1).Value = 7
$ws.Cells.Item(313, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(313, 3).Value = "Ñuble"
$ws.Cells.Item(313, 4).Value = 44461
$ws.Cells.Item(313, 5).Value = 16
$ws.Cells.Item(313, 6).Value = 100112004
$ws.Cells.Item(313, 7).Value = "Cebolla"
$ws.Cells.Item(313, 8).Value = "Sin especificar"
$ws.Cells.Item(313, 9).Value = "1a (guarda)"
$ws.Cells.Item(313, 10).Value = 300
$ws.Cells.Item(313, 11).Value = 3000
$ws.Cells.Item(313, 12).Value = 3200
$ws.Cells.Item(313, 13).Value = 3100
$ws.Cells.Item(313, 14).Value = "`$/malla 15 kilos"
$ws.Cells.Item(313, 15).Value = "Región del Maule"
$ws.Cells.Item(313, 16).Value = 207
$ws.Cells.Item(313, 17).Value = 15
$ws.Cells.Item(313, 18).Value = "Hortaliza"
